# Flight_Mission_Cycle.xlsx edit: replace the "Settings" sheet with a new
# "Writing" sheet holding Force-Time graph data, and make it the active tab.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove the old "Settings" sheet -----------------------------------
$settings = $wb.Worksheets.Item("Settings")
[void]$settings.Delete()

# --- Add the new "Writing" sheet at the end -----------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$writing = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$writing.Name = "Writing"

# --- Populate the Force-Time graph data ---------------------------------
$writing.Range("A1").Value = "Writing"

$writing.Range("A2").Value = "Force_End"
$writing.Range("B2").Value = 50
$writing.Range("C2").Value = 50
$writing.Range("D2").Value = 100
$writing.Range("E2").Value = 0

$writing.Range("A3").Value = "Duration"
$writing.Range("B3").Value = 20
$writing.Range("C3").Value = 30
$writing.Range("D3").Value = 40
$writing.Range("E3").Value = 25

# Match the authored column width of the "Settings" sheet (OOXML width 11)
$writing.Columns.Item(1).ColumnWidth = 10.1667

# --- Make "Writing" the active/selected sheet ---------------------------
[void]$writing.Activate()
[void]$writing.Range("F14").Select()
